$d = $word.ActiveDocument

$replacements = @(
    @{old = "255÷7=36, 3"; new = "442÷3=147, 1"},
    @{old = "287÷8=35, 7"; new = "826÷8=103, 2"},
    @{old = "487÷7=69, 4"; new = "953÷3=317, 2"},
    @{old = "507÷5=101, 2"; new = "817÷3=272, 1"},
    @{old = "445÷3=148, 1"; new = "518÷7=74, 0"},
    @{old = "475÷5=95, 0"; new = "960÷7=137, 1"},
    @{old = "167÷6=27, 5"; new = "390÷3=130, 0"},
    @{old = "427÷2=213, 1"; new = "961÷2=480, 1"},
    @{old = "978÷2=489, 0"; new = "780÷8=97, 4"},
    @{old = "517÷7=73, 6"; new = "988÷9=109, 7"},
    @{old = "283÷8=35, 3"; new = "469÷8=58, 5"},
    @{old = "158÷5=31, 3"; new = "394÷3=131, 1"},
    @{old = "322÷3=107, 1"; new = "212÷9=23, 5"},
    @{old = "763÷2=381, 1"; new = "983÷9=109, 2"},
    @{old = "768÷7=109, 5"; new = "942÷5=188, 2"},
    @{old = "976÷4=244, 0"; new = "473÷8=59, 1"},
    @{old = "914÷6=152, 2"; new = "167÷4=41, 3"},
    @{old = "930÷9=103, 3"; new = "185÷9=20, 5"},
    @{old = "198÷3=66, 0"; new = "805÷6=134, 1"},
    @{old = "702÷5=140, 2"; new = "503÷5=100, 3"},
    @{old = "430÷8=53, 6"; new = "623÷4=155, 3"},
    @{old = "702÷8=87, 6"; new = "197÷9=21, 8"},
    @{old = "418÷7=59, 5"; new = "156÷4=39, 0"},
    @{old = "869÷3=289, 2"; new = "351÷9=39, 0"},
    @{old = "663÷5=132, 3"; new = "343÷7=49, 0"}
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
